$wb = $excel.ActiveWorkbook

# The new sheet ("FTNC_Average_Demand56") is structurally identical to the
# last existing sheet ("FTNC_Average_Demand55" - same headers, same style,
# same page setup), just with a fresh set of numbers, so clone that sheet
# (preserves sheetPr/styles/pageMargins exactly) and place the copy right
# after it, then overwrite its name + data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FTNC_Average_Demand56"

# Row 1 headers / row 2 label (A2) are already correct, carried over from the
# copied sheet. Only the numeric results in row 2 change.
$newSheet.Range("B2").Value = 2289.585598990748
$newSheet.Range("C2").Value = 12867.73307817783
$newSheet.Range("D2").Value = 624.6448382569154
$newSheet.Range("E2").Value = 11.68808730780301
$newSheet.Range("F2").Value = 15793.65160273333
